$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.443.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.700.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "693.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.699.11"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.43"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.55"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.323.71"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.695.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.505.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.113"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.29"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "481.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.02"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.665"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.71%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.848.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.99%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000131"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.43"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -8.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -10.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.73"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -9.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.86"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.08"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.10"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.166"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.668.66"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.50"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.41"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.33"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0933"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.954"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "164.10"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "30.21"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.82"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -14.75%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.35"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.30%  "

$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000286"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -8.03%  "

Write-Host "Applied cryptos update"